$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: a new "Alcachofa" price entry for Macroferia Regional de
# Talca is inserted right after the existing row 22, pushing every
# following row down by one. The new row starts out as a duplicate of the
# row that is about to land below it (old row 23), and only its date
# (column D) is updated to the new reporting date.

$ws.Rows(23).Insert()
$ws.Range("A24:R24").Copy($ws.Range("A23:R23"))

# 2021-09-13 -> Excel serial date 44452
$ws.Range("D23").Value2 = 44452
